$wb = $excel.ActiveWorkbook
$wsInputs = $wb.Worksheets.Item("Inputs & Income Statement")

# --- Update the two input assumptions (buying/selling price base year) ---
$wsInputs.Range("C5").Value = 18
$wsInputs.Range("C6").Value = 42

# --- Rewrite the Question rich-text cell (B18) with new wording/figures ---
$cell = $wsInputs.Range("B18")
$cell.Value = 'Question: The selling price of 1 kg of rice is 42-220 Rs(min-max set accordingly)and the buying price is 18-26 Rs.(min-max set accordingly)  In the first year, the ABC company sells 2000 tonnes. Additionally, there is a Maintanace Charge(Storage & Labour) of 2% and a Delivery Charge(Transpotation & Labour ) of 3%  of buying price that reduce from the total profit. If historical data shows that every year there is an 11% increment in the buying price, so the ABC company also increases the selling price accordingly. Apart from that, there is an increment of 0.5% in the Maintanace charge and 0.8% in the delivery charge each year. How much revenue can the company make if in the second year it sells 3000 tonnes, in the third year it sells 2500 tonnes, in the fourth year it sells 3100 tonnes, and in the fifth year it sells 1800 tonnes? if each year tax is 4 %'
$cell.Characters(1,8).Font.Bold = $true
$cell.Characters(1,8).Font.Size = 16
$cell.Characters(9,2).Font.Bold = $true
$cell.Characters(9,2).Font.Size = 12
$cell.Characters(11,21).Font.Bold = $false
$cell.Characters(11,21).Font.Size = 14
$cell.Characters(32,6).Font.Bold = $true
$cell.Characters(32,6).Font.Size = 14
$cell.Characters(38,10).Font.Bold = $false
$cell.Characters(38,10).Font.Size = 14
$cell.Characters(48,34).Font.Bold = $true
$cell.Characters(48,34).Font.Size = 14
$cell.Characters(82,23).Font.Bold = $false
$cell.Characters(82,23).Font.Size = 14
$cell.Characters(105,36).Font.Bold = $true
$cell.Characters(105,36).Font.Size = 14
$cell.Characters(141,24).Font.Bold = $false
$cell.Characters(141,24).Font.Size = 14
$cell.Characters(165,3).Font.Bold = $true
$cell.Characters(165,3).Font.Size = 14
$cell.Characters(168,15).Font.Bold = $false
$cell.Characters(168,15).Font.Size = 14
$cell.Characters(183,4).Font.Bold = $true
$cell.Characters(183,4).Font.Size = 14
$cell.Characters(187,1).Font.Bold = $false
$cell.Characters(187,1).Font.Size = 14
$cell.Characters(188,6).Font.Bold = $true
$cell.Characters(188,6).Font.Size = 14
$cell.Characters(194,27).Font.Bold = $false
$cell.Characters(194,27).Font.Size = 14
$cell.Characters(221,36).Font.Bold = $true
$cell.Characters(221,36).Font.Size = 14
$cell.Characters(257,3).Font.Bold = $false
$cell.Characters(257,3).Font.Size = 14
$cell.Characters(260,2).Font.Bold = $true
$cell.Characters(260,2).Font.Size = 14
$cell.Characters(262,6).Font.Bold = $false
$cell.Characters(262,6).Font.Size = 14
$cell.Characters(268,42).Font.Bold = $true
$cell.Characters(268,42).Font.Size = 14
$cell.Characters(310,3).Font.Bold = $false
$cell.Characters(310,3).Font.Size = 14
$cell.Characters(313,2).Font.Bold = $true
$cell.Characters(313,2).Font.Size = 14
$cell.Characters(315,106).Font.Bold = $false
$cell.Characters(315,106).Font.Size = 14
$cell.Characters(421,4).Font.Bold = $true
$cell.Characters(421,4).Font.Size = 14
$cell.Characters(425,38).Font.Bold = $false
$cell.Characters(425,38).Font.Size = 14
$cell.Characters(463,3).Font.Bold = $true
$cell.Characters(463,3).Font.Size = 14
$cell.Characters(466,97).Font.Bold = $false
$cell.Characters(466,97).Font.Size = 14
$cell.Characters(563,4).Font.Bold = $true
$cell.Characters(563,4).Font.Size = 14
$cell.Characters(567,30).Font.Bold = $false
$cell.Characters(567,30).Font.Size = 14
$cell.Characters(597,4).Font.Bold = $true
$cell.Characters(597,4).Font.Size = 14
$cell.Characters(601,104).Font.Bold = $false
$cell.Characters(601,104).Font.Size = 14
$cell.Characters(705,4).Font.Bold = $true
$cell.Characters(705,4).Font.Size = 14
$cell.Characters(709,1).Font.Bold = $false
$cell.Characters(709,1).Font.Size = 14
$cell.Characters(710,6).Font.Bold = $true
$cell.Characters(710,6).Font.Size = 14
$cell.Characters(716,29).Font.Bold = $false
$cell.Characters(716,29).Font.Size = 14
$cell.Characters(745,4).Font.Bold = $true
$cell.Characters(745,4).Font.Size = 14
$cell.Characters(749,1).Font.Bold = $false
$cell.Characters(749,1).Font.Size = 14
$cell.Characters(750,6).Font.Bold = $true
$cell.Characters(750,6).Font.Size = 14
$cell.Characters(756,30).Font.Bold = $false
$cell.Characters(756,30).Font.Size = 14
$cell.Characters(786,4).Font.Bold = $true
$cell.Characters(786,4).Font.Size = 14
$cell.Characters(790,1).Font.Bold = $false
$cell.Characters(790,1).Font.Size = 14
$cell.Characters(791,6).Font.Bold = $true
$cell.Characters(791,6).Font.Size = 14
$cell.Characters(797,33).Font.Bold = $false
$cell.Characters(797,33).Font.Size = 14
$cell.Characters(830,4).Font.Bold = $true
$cell.Characters(830,4).Font.Size = 14
$cell.Characters(834,1).Font.Bold = $false
$cell.Characters(834,1).Font.Size = 14
$cell.Characters(835,6).Font.Bold = $true
$cell.Characters(835,6).Font.Size = 14
$cell.Characters(841,1).Font.Bold = $false
$cell.Characters(841,1).Font.Size = 14
$cell.Characters(842,24).Font.Bold = $true
$cell.Characters(842,24).Font.Size = 14

# --- Update the view: scroll/select to show the Question block ---
$wsInputs.Activate()
$wsInputs.Range("B18:H28").Select()
